$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 5048
$ws.Range("E2").Value = 103
$ws.Range("F2").Value = 103
$ws.Range("G2").Value = 88
$ws.Range("H2").Value = 26
$ws.Range("I2").Value = 96
$ws.Range("J2").Value = -71
$ws.Range("K2").Value = 5216
$ws.Range("L2").Value = 3230
$ws.Range("M2").Value = 1986
$ws.Range("N2").Value = 1978
$ws.Range("O2").Value = 8
$ws.Range("P2").Value = 413
$ws.Range("Q2").Value = 285
$ws.Range("R2").Value = -640
$ws.Range("S2").Value = 503
$ws.Range("T2").Value = 767
$ws.Range("U2").Value = -482
$ws.Range("V2").Value = 2104
$ws.Range("W2").Value = 2.04
$ws.Range("X2").Value = 0.51
$ws.Range("Y2").Value = 4.84
$ws.Range("Z2").Value = 0.52
$ws.Range("AA2").Value = 162.61
$ws.Range("AB2").Value = 392.34
$ws.Range("AC2").Value = 233
$ws.Range("AD2").Value = 20.41
$ws.Range("AE2").Value = 4965
$ws.Range("AF2").Value = 0.96
$ws.Range("AG2").Value = 80
$ws.Range("AH2").Value = 1.68
$ws.Range("AI2").Value = 33.16
$ws.Range("AJ2").Value = 41268398
$ws.Range("D3").Value = 5209
$ws.Range("E3").Value = 202
$ws.Range("F3").Value = 202
$ws.Range("G3").Value = 159
$ws.Range("H3").Value = 65
$ws.Range("I3").Value = 156
$ws.Range("J3").Value = -92
$ws.Range("K3").Value = 5063
$ws.Range("L3").Value = 3082
$ws.Range("M3").Value = 1981
$ws.Range("N3").Value = 1967
$ws.Range("O3").Value = 14
$ws.Range("P3").Value = 413
$ws.Range("Q3").Value = 488
$ws.Range("R3").Value = -312
$ws.Range("S3").Value = -132
$ws.Range("T3").Value = 291
$ws.Range("U3").Value = 198
$ws.Range("V3").Value = 2030
$ws.Range("W3").Value = 3.87
$ws.Range("X3").Value = 1.24
$ws.Range("Y3").Value = 7.92
$ws.Range("Z3").Value = 1.26
$ws.Range("AA3").Value = 155.55
$ws.Range("AB3").Value = 409.48
$ws.Range("AC3").Value = 379
$ws.Range("AD3").Value = 12.19
$ws.Range("AE3").Value = 4937
$ws.Range("AF3").Value = 0.93
$ws.Range("AG3").Value = 90
$ws.Range("AH3").Value = 1.95
$ws.Range("AI3").Value = 22.94
$ws.Range("AJ3").Value = 41268398
$ws.Range("D4").Value = 5542
$ws.Range("E4").Value = 111
$ws.Range("F4").Value = 111
$ws.Range("G4").Value = 70
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 37
$ws.Range("J4").Value = -35
$ws.Range("K4").Value = 5511
$ws.Range("L4").Value = 3530
$ws.Range("M4").Value = 1981
$ws.Range("N4").Value = 1976
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 413
$ws.Range("Q4").Value = -1
$ws.Range("R4").Value = -105
$ws.Range("S4").Value = 114
$ws.Range("T4").Value = 190
$ws.Range("U4").Value = -191
$ws.Range("V4").Value = 2245
$ws.Range("W4").Value = 1.99
$ws.Range("X4").Value = 0.04
$ws.Range("Y4").Value = 1.89
$ws.Range("Z4").Value = 0.05
$ws.Range("AA4").Value = 178.24
$ws.Range("AB4").Value = 396.74
$ws.Range("AC4").Value = 90
$ws.Range("AD4").Value = 47.03
$ws.Range("AE4").Value = 4970
$ws.Range("AF4").Value = 0.85
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 2.36
$ws.Range("AI4").Value = 106.89
$ws.Range("AJ4").Value = 41268398
$ws.Range("D5").Value = 5313
$ws.Range("E5").Value = -80
$ws.Range("F5").Value = -80
$ws.Range("G5").Value = -118
$ws.Range("H5").Value = -149
$ws.Range("I5").Value = -140
$ws.Range("J5").Value = -8
$ws.Range("K5").Value = 4833
$ws.Range("L5").Value = 3102
$ws.Range("M5").Value = 1731
$ws.Range("N5").Value = 1727
$ws.Range("O5").Value = 5
$ws.Range("P5").Value = 413
$ws.Range("Q5").Value = 67
$ws.Range("R5").Value = 110
$ws.Range("S5").Value = -90
$ws.Range("T5").Value = 391
$ws.Range("U5").Value = -324
$ws.Range("V5").Value = 2154
$ws.Range("W5").Value = -1.51
$ws.Range("X5").Value = -2.8
$ws.Range("Y5").Value = -7.59
$ws.Range("Z5").Value = -2.88
$ws.Range("AA5").Value = 179.16
$ws.Range("AB5").Value = 353.35
$ws.Range("AC5").Value = -340
$ws.Range("AD5").Value = -12.41
$ws.Range("AE5").Value = 4383
$ws.Range("AF5").Value = 0.96
$ws.Range("AG5").Value = 70
$ws.Range("AH5").Value = 1.66
$ws.Range("AI5").Value = -19.63
$ws.Range("AJ5").Value = 41268398
$ws.Range("D6").Value = 5603
$ws.Range("E6").Value = 75
$ws.Range("F6").Value = 75
$ws.Range("G6").Value = -1
$ws.Range("H6").Value = -74
$ws.Range("I6").Value = -68
$ws.Range("K6").Value = 5056
$ws.Range("L6").Value = 3279
$ws.Range("M6").Value = 1777
$ws.Range("N6").Value = 1776
$ws.Range("P6").Value = 413
$ws.Range("Q6").Value = -266
$ws.Range("R6").Value = -167
$ws.Range("S6").Value = 15
$ws.Range("T6").Value = 258
$ws.Range("U6").Value = -524
$ws.Range("V6").Value = 2188
$ws.Range("W6").Value = 1.34
$ws.Range("X6").Value = -1.32
$ws.Range("Y6").Value = -3.86
$ws.Range("Z6").Value = -1.49
$ws.Range("AA6").Value = 184.59
$ws.Range("AB6").Value = 320.72
$ws.Range("AC6").Value = -164
$ws.Range("AD6").Value = -39.69
$ws.Range("AE6").Value = 4466
$ws.Range("AF6").Value = 1.46
$ws.Range("AG6").Value = 70
$ws.Range("AH6").Value = 1.08
$ws.Range("AI6").Value = -41.12
$ws.Range("AJ6").Value = 41268398

# Rows 7-9: clear all data columns (D:AJ), keeping only A/B/C (row#, category labels)
$ws.Range("D7:AJ9").ClearContents()
